$d = $word.ActiveDocument

$d.Content.Find.Execute("809×7=5663", $true, $false, $false, $false, $false, $true, 1, $false, "304×3=912", 2) | Out-Null
$d.Content.Find.Execute("516×2=1032", $true, $false, $false, $false, $false, $true, 1, $false, "434×5=2170", 2) | Out-Null
$d.Content.Find.Execute("174×7=1218", $true, $false, $false, $false, $false, $true, 1, $false, "343×4=1372", 2) | Out-Null
$d.Content.Find.Execute("419×3=1257", $true, $false, $false, $false, $false, $true, 1, $false, "358×9=3222", 2) | Out-Null
$d.Content.Find.Execute("799×9=7191", $true, $false, $false, $false, $false, $true, 1, $false, "849×8=6792", 2) | Out-Null
$d.Content.Find.Execute("548×7=3836", $true, $false, $false, $false, $false, $true, 1, $false, "591×9=5319", 2) | Out-Null
$d.Content.Find.Execute("188×7=1316", $true, $false, $false, $false, $false, $true, 1, $false, "899×6=5394", 2) | Out-Null
$d.Content.Find.Execute("352×3=1056", $true, $false, $false, $false, $false, $true, 1, $false, "207×8=1656", 2) | Out-Null
$d.Content.Find.Execute("543×5=2715", $true, $false, $false, $false, $false, $true, 1, $false, "333×8=2664", 2) | Out-Null
$d.Content.Find.Execute("222×6=1332", $true, $false, $false, $false, $false, $true, 1, $false, "744×5=3720", 2) | Out-Null
$d.Content.Find.Execute("422×6=2532", $true, $false, $false, $false, $false, $true, 1, $false, "862×8=6896", 2) | Out-Null
$d.Content.Find.Execute("736×5=3680", $true, $false, $false, $false, $false, $true, 1, $false, "558×2=1116", 2) | Out-Null
$d.Content.Find.Execute("691×9=6219", $true, $false, $false, $false, $false, $true, 1, $false, "492×9=4428", 2) | Out-Null
$d.Content.Find.Execute("439×8=3512", $true, $false, $false, $false, $false, $true, 1, $false, "243×9=2187", 2) | Out-Null
$d.Content.Find.Execute("161×9=1449", $true, $false, $false, $false, $false, $true, 1, $false, "940×3=2820", 2) | Out-Null
$d.Content.Find.Execute("850×2=1700", $true, $false, $false, $false, $false, $true, 1, $false, "440×3=1320", 2) | Out-Null
$d.Content.Find.Execute("297×9=2673", $true, $false, $false, $false, $false, $true, 1, $false, "516×7=3612", 2) | Out-Null
$d.Content.Find.Execute("929×9=8361", $true, $false, $false, $false, $false, $true, 1, $false, "236×5=1180", 2) | Out-Null
$d.Content.Find.Execute("147×3=441", $true, $false, $false, $false, $false, $true, 1, $false, "426×8=3408", 2) | Out-Null
$d.Content.Find.Execute("914×6=5484", $true, $false, $false, $false, $false, $true, 1, $false, "841×8=6728", 2) | Out-Null
$d.Content.Find.Execute("596×8=4768", $true, $false, $false, $false, $false, $true, 1, $false, "745×3=2235", 2) | Out-Null
$d.Content.Find.Execute("139×5=695", $true, $false, $false, $false, $false, $true, 1, $false, "330×6=1980", 2) | Out-Null
$d.Content.Find.Execute("750×3=2250", $true, $false, $false, $false, $false, $true, 1, $false, "791×6=4746", 2) | Out-Null
$d.Content.Find.Execute("180×7=1260", $true, $false, $false, $false, $false, $true, 1, $false, "650×4=2600", 2) | Out-Null
$d.Content.Find.Execute("520×6=3120", $true, $false, $false, $false, $false, $true, 1, $false, "502×9=4518", 2) | Out-Null

Write-Host "Done replacing all values."
